# Updated cryptos list with refreshed price/volume data (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "48.210.86"
$ws.Cells.Item(2, 5).Value = "  +1.93%  "
$ws.Cells.Item(3, 4).Value = "2.530.81"
$ws.Cells.Item(3, 5).Value = "  +1.03%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.08%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "323.86"
$ws.Cells.Item(5, 5).Value = "  -0.13%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "109.28"
$ws.Cells.Item(6, 5).Value = "  +0.06%  "
$ws.Cells.Item(7, 5).Value = "  +0.48%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.999"
$ws.Cells.Item(8, 5).Value = "  -0.05%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.565"
$ws.Cells.Item(9, 5).Value = "  +5.03%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "40.53"
$ws.Cells.Item(10, 5).Value = "  +3.55%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "20.28"
$ws.Cells.Item(11, 5).Value = "  +10.13%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0823"
$ws.Cells.Item(12, 5).Value = "  +1.18%  "
$ws.Cells.Item(13, 5).Value = "  +1.10%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "7.29"
$ws.Cells.Item(14, 5).Value = "  +1.28%  "
$ws.Cells.Item(15, 4).Value = "2.921.65"
$ws.Cells.Item(15, 5).Value = "  +0.83%  "
$ws.Cells.Item(16, 4).Value = "2.534.39"
$ws.Cells.Item(16, 5).Value = "  +0.97%  "
$ws.Cells.Item(17, 5).Value = "  +0.85%  "
$ws.Cells.Item(18, 4).Value = "48.134.84"
$ws.Cells.Item(18, 5).Value = "  +1.92%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "13.28"
$ws.Cells.Item(19, 5).Value = "  +3.31%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "6.64"
$ws.Cells.Item(20, 5).Value = "  -0.19%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0951"
$ws.Cells.Item(21, 5).Value = "  +0.84%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "2.73"
$ws.Cells.Item(22, 5).Value = "  +0.44%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "72.48"
$ws.Cells.Item(23, 5).Value = "  +2.80%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "270.87"
$ws.Cells.Item(24, 5).Value = "  +9.46%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.59"
$ws.Cells.Item(25, 5).Value = "  -0.33%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "26.27"
$ws.Cells.Item(26, 5).Value = "  +0.95%  "
$ws.Cells.Item(27, 5).Value = "  +0.08%  "

# Rows 28-31: coin entries swapped positions with updated D/E values
$ws.Cells.Item(28, 2).Value = "Cosmos"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "10.17"
$ws.Cells.Item(28, 5).Value = "  +1.22%  "

$ws.Cells.Item(29, 2).Value = "Kaspa"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.148"
$ws.Cells.Item(29, 5).Value = "  +7.14%  "

$ws.Cells.Item(30, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "35.33"
$ws.Cells.Item(30, 5).Value = "  -0.73%  "

$ws.Cells.Item(31, 2).Value = "Toncoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "2.10"
$ws.Cells.Item(31, 5).Value = "  -8.74%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "49.77"
$ws.Cells.Item(32, 5).Value = "  -0.20%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "20.02"
$ws.Cells.Item(33, 5).Value = "  -0.16%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.41"
$ws.Cells.Item(34, 5).Value = "  -0.72%  "
$ws.Cells.Item(35, 5).Value = "  -0.16%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.0793"
$ws.Cells.Item(36, 5).Value = "  +0.02%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.00"
$ws.Cells.Item(37, 5).Value = "  +0.58%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "4.75"
$ws.Cells.Item(38, 5).Value = "  +0.84%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "3.02"
$ws.Cells.Item(39, 5).Value = "  +0.43%  "
$ws.Cells.Item(40, 5).Value = "  +0.08%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "22.63"
$ws.Cells.Item(41, 5).Value = "  +6.39%  "
$ws.Cells.Item(42, 5).Value = "  -2.19%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "117.91"
$ws.Cells.Item(43, 5).Value = "  -2.85%  "
$ws.Cells.Item(44, 5).Value = "  +0.20%  "
$ws.Cells.Item(45, 4).Value = "2.010.96"
$ws.Cells.Item(45, 5).Value = "  +0.71%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "3.17"
$ws.Cells.Item(46, 5).Value = "  +2.58%  "
$ws.Cells.Item(47, 5).Value = "  +6.87%  "
$ws.Cells.Item(48, 5).Value = "  -1.22%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "9.12"
$ws.Cells.Item(49, 5).Value = "  +0.34%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "5.25"
$ws.Cells.Item(50, 5).Value = "  +0.57%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "80.45"
$ws.Cells.Item(51, 5).Value = "  +3.18%  "
